$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45044
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("S2").Value = 833

# Row 3
$ws.Range("D3").Value = 45043
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("S3").Value = 833

# Row 4
$ws.Range("D4").Value = 45030
$ws.Range("M4").Value = 40
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 18000
$ws.Range("S4").Value = 1000

# Row 5
$ws.Range("D5").Value = 45041
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("S5").Value = 833

# Row 6
$ws.Range("D6").Value = 45001

# Row 7
$ws.Range("D7").Value = 45028
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 18000
$ws.Range("O7").Value = 18000
$ws.Range("P7").Value = 18000
$ws.Range("S7").Value = 1000

# Row 8
$ws.Range("D8").Value = 45014
$ws.Range("M8").Value = 30
$ws.Range("N8").Value = 18000
$ws.Range("O8").Value = 18000
$ws.Range("P8").Value = 18000
$ws.Range("Q8").Value = "$/caja 18 kilos"
$ws.Range("R8").Value = "Región Metropolitana"
$ws.Range("S8").Value = 1000
$ws.Range("T8").Value = 18

# Row 10
$ws.Range("D10").Value = 45002
$ws.Range("M10").Value = 30

# Row 11
$ws.Range("D11").Value = 45049
$ws.Range("M11").Value = 80

# Row 12
$ws.Range("D12").Value = 44999
$ws.Range("M12").Value = 60
$ws.Range("N12").Value = 17000
$ws.Range("P12").Value = 17500
$ws.Range("S12").Value = 972

# Row 13
$ws.Range("D13").Value = 45020
$ws.Range("M13").Value = 50
$ws.Range("Q13").Value = "$/caja 16 kilos"
$ws.Range("R13").Value = "Provincia de Los Andes"
$ws.Range("S13").Value = 938
$ws.Range("T13").Value = 16

# Row 14
$ws.Range("D14").Value = 45062
$ws.Range("M14").Value = 90
$ws.Range("N14").Value = 13000
$ws.Range("O14").Value = 14000
$ws.Range("P14").Value = 13444
$ws.Range("S14").Value = 747

# Row 15
$ws.Range("D15").Value = 45033

# Row 16
$ws.Range("D16").Value = 45036
$ws.Range("R16").Value = "Región Metropolitana"

# Row 17
$ws.Range("D17").Value = 45021
$ws.Range("M17").Value = 60
$ws.Range("N17").Value = 15000
$ws.Range("O17").Value = 16000
$ws.Range("P17").Value = 15500
$ws.Range("R17").Value = "Provincia de Los Andes"
$ws.Range("S17").Value = 861

# Row 18
$ws.Range("D18").Value = 45037
$ws.Range("N18").Value = 16000
$ws.Range("O18").Value = 16000
$ws.Range("P18").Value = 16000
$ws.Range("S18").Value = 889
